$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Marking" row: correct-answer mark value (B11)
$ws.Range("B11").Value = 5

# "Total" row: total marks obtained (B12) and the "obtained/max" label (E12)
$ws.Range("B12").Value = 140
$ws.Range("E12").Value = "140/140"
